# Datacamp/disk_savvy courses tracker update
#
# The diff inserts a brand-new row at row 48 (pushing the previous rows
# 48-54 down to 49-55) and fills the new row with a new course entry:
#   A48 = "Exploratory Data Analysis in Python"  (a brand-new shared string)
#   K48 = 2
# The row that used to be row 48 (an empty, style-3 "A" cell acting as a
# section separator) ends up as the new row 49, unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 48. Excel copies the
# formatting of the row above (row 47, style 3 on column A / style 2 on
# column K) down onto both the new row and the row that gets pushed down,
# which mirrors the style pattern seen in the target file.
$ws.Rows("48:48").Insert()

# Fill in the freshly inserted row with the new course/data.
$ws.Range("A48").Value2 = "Exploratory Data Analysis in Python"
$ws.Range("K48").Value2 = 2

# Reflect the updated view state: the sheet had scrolled down and the
# user's selection moved to A49 after the insertion.
$excel.ActiveWindow.ScrollRow = 41
$ws.Range("A49").Select()
